$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3: the "Ticketart" dropdown option list changes its separator
# from commas to slashes (Bus,Zug,U-Bahn -> Bus/Zug/U-Bahn), while the
# "Häufigkeit" dropdown part stays unchanged.
$ws.Range("D3").Value = "Ticketart:dropdown(Bus/Zug/U-Bahn);Häufigkeit:dropdown(Täglich,Wöchentlich,Selten)"

# Reflect the new active selection/cell recorded in the sheet view.
$ws.Range("D3").Select()
